# Daily attendance processing - 2025-12-31 18:40:49
# Reorders the "Recorded By" (column G) list of names/emails for each row:
#   - if "System" (exact case) appears in the list, move it to the end,
#     keeping the relative order of the remaining entries
#   - otherwise (e.g. two emails with no "System"), reverse the order

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G ("Recorded By")
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notmatch ",") { continue }

    $parts = @($val -split ", ")

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }

    $newParts = @()
    if ($hasSystem) {
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) { $newParts += $p }
        }
        foreach ($p in $parts) {
            if ($p.Equals("System")) { $newParts += $p }
        }
    } else {
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newVal = [string]::Join(", ", $newParts)
    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
